# Updates D (Price) and E (Volume(1h)) columns on Sheet1 for the
# refreshed "cryptos" list (GitHub Actions data refresh).
#
# Prices such as "67.641.19" or "0.143" look numeric to Excel, so a
# plain .Value assignment would get silently parsed/rounded as a
# number. We force the destination cell to Text format first, write
# the literal string, then ClearFormats() so the cell ends up with
# no explicit style (matching the original un-styled data cells)
# while keeping the text exactly as authored.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.641.19'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.11%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.515.06'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.25%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '609.86'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.57'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.82%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.513.41'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.27%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  +1.07%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.143'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +2.78%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.62'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +7.89%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.431'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.43%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '32.73'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +2.37%  '
$ws.Range("E14").Value = '  -1.48%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.108.64'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.42%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.517.62'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.01%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.620.38'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.16%  '
$ws.Range("E18").Value = '  +0.47%  '
$ws.Range("E19").Value = '  +2.40%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.60'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +2.52%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.89'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +6.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '449.11'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.58%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.634'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.64%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.28'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +1.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000128'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.75%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.653.99'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.31%  '
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.88'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +6.34%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.13'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.60%  '
$ws.Range("E30").Value = '  +0.86%  '
$ws.Range("E31").Value = '  +5.82%  '
$ws.Range("E32").Value = '  +2.73%  '
$ws.Range("E33").Value = '  -0.04%  '
$ws.Range("E34").Value = '  +0.11%  '
$ws.Range("E35").Value = '  +1.04%  '
$ws.Range("E36").Value = '  +1.69%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.508.51'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.27%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.06'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.12%  '
$ws.Range("E40").Value = '  +5.48%  '
$ws.Range("E41").Value = '  +0.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0899'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +2.69%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '174.25'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.67%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.50'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +1.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '30.37'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +11.66%  '
$ws.Range("E46").Value = '  +0.37%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '46.74'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.59%  '
$ws.Range("E48").Value = '  +3.17%  '
$ws.Range("E49").Value = '  -3.23%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.68'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.19%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.254'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.98%  '
